# Applies the "Updates to Jun's files 10.21" revision to
# Elast of Component E Demand wrt E Cost.xlsx

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("B4").Value = 2014
$about.Range("B6").Value = "https://www.eia.gov/analysis/studies/buildings/energyuse/pdf/price_elasticities.pdf"
$about.Range("B5").Value = "Price Elasticities for Energy Use in Buildings of the United States"
$about.Range("B7").Value = "Appendix"

$about.Range("A10").Value = "We use same-price, long-run elasticities minus the 3-year short-run elasticities."
$about.Range("A11").Value = "We calculate it this way because we assume that 3-year elasticities primarily reflect behavior"
$about.Range("A14").Value = "all timescales.  So, the portion of the long-run elasticitiy represented by the 3-year elasticity"

# The new citation text on rows 4-7 is a different length than before, so Excel
# re-autofits those row heights (dropping the explicit "ht" override).
$about.Range("A4:A7").EntireRow.AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "EIA Table 1"
# ---------------------------------------------------------------------------
$eia = $wb.Worksheets.Item("EIA Table 1")

# Residential block
$eia.Range("B7").Value = -0.12
$eia.Range("C7").Value = -0.21
$eia.Range("D7").Value = -0.25
$eia.Range("E7").Value = -0.28
$eia.Range("F7").Value = 0
$eia.Range("G7").Value = 0

$eia.Range("B8").Value = -0.07
$eia.Range("C8").Value = -0.13
$eia.Range("D8").Value = -0.15
$eia.Range("E8").Value = 0.03
$eia.Range("F8").Value = -0.21
$eia.Range("G8").Value = 0

$eia.Range("B9").Value = -0.07
$eia.Range("C9").Value = -0.12
$eia.Range("D9").Value = -0.14
$eia.Range("E9").Value = 0
$eia.Range("F9").Value = 0
$eia.Range("G9").Value = -0.22

# Commercial block
$eia.Range("B14").Value = -0.11
$eia.Range("C14").Value = -0.18
$eia.Range("D14").Value = -0.22
$eia.Range("E14").Value = -0.33
$eia.Range("F14").Value = 0.09
$eia.Range("G14").Value = 0

$eia.Range("B15").Value = -0.15
$eia.Range("C15").Value = -0.25
$eia.Range("D15").Value = -0.3
$eia.Range("E15").Value = 0.15
$eia.Range("F15").Value = -0.58
$eia.Range("G15").Value = 0.02

$eia.Range("B16").Value = -0.14
$eia.Range("C16").Value = -0.24
$eia.Range("D16").Value = -0.29
$eia.Range("E16").Value = 0
$eia.Range("F16").Value = 0.05
$eia.Range("G16").Value = -0.42

# ---------------------------------------------------------------------------
# Sheet "EoCEDwEC"
# ---------------------------------------------------------------------------
$eoc = $wb.Worksheets.Item("EoCEDwEC")

$eoc.Range("B2").Formula = "='EIA Table 1'!E7-'EIA Table 1'!D7"
$eoc.Range("D2").Formula = "='EIA Table 1'!E14-'EIA Table 1'!D14"
$eoc.Range("B4").Formula = "='EIA Table 1'!F8-'EIA Table 1'!D8"
$eoc.Range("D4").Formula = "='EIA Table 1'!F15-'EIA Table 1'!D15"
$eoc.Range("B5").Formula = "='EIA Table 1'!G9-'EIA Table 1'!D9"
$eoc.Range("D5").Formula = "='EIA Table 1'!G16-'EIA Table 1'!D16"

# ---------------------------------------------------------------------------
# View state: EoCEDwEC becomes the active/selected sheet, with a couple of
# incidental cell selections recorded on the other two sheets.
# ---------------------------------------------------------------------------
$about.Range("A27").Select() | Out-Null
$eia.Range("E17").Select() | Out-Null
$eoc.Activate()
$eoc.Range("H29").Select() | Out-Null
